$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I1 = "I0" and J1 = "IF", copying the existing header
# (H1) formatting so the new cells get the same bold/bordered/centered style.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "IF"

# Data rows 2-28: I column is always 1, J column mirrors the H column (IP).
for ($r = 2; $r -le 28; $r++) {
    $hval = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hval
}
